$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: record a new comparison row (h5) and tweak F5 ---
$ws2.Range("F5").Value = "Data Ds"

$ws2.Range("A8").Value = "h5"
$ws2.Range("B8").Value = ""
$ws2.Range("A8:B8").HorizontalAlignment = -4108
$ws2.Range("C8").Value = 2
$ws2.Range("D8").Value = 3
$ws2.Range("E8").Value = 4
$ws2.Range("F8").Value = 5
$ws2.Range("A8:B8").Merge()

# --- Selections / active sheet ---
$ws1.Range("L6").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 181
$ws2.Range("F8").Select()
